# Daily update at 8 AM UTC
# Adds the next day's row (row 52) to the "Wins Over Time" sheet and
# normalizes the previous last row's date style back to the standard one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 was previously the last row and used the "last row" date style
# (plain YYYY-MM-DD). Since it is no longer the last row, restore it to the
# standard date/time style used by all the other data rows.
$ws.Range("A51").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 52.
$ws.Range("A52").Value = 45792
$ws.Range("B52").Value = 215
$ws.Range("C52").Value = 218
$ws.Range("D52").Value = 221

# The new last row gets the distinct "last row" date style.
$ws.Range("A52").NumberFormat = "YYYY-MM-DD"
